$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: change date separators from "/" to "-" (28/07/2022 -> 28-07-2022) ---
#
# Rows whose day-of-month is > 12 are unambiguous, so a plain text Replace
# leaves them as text. Rows whose day-of-month is <= 12 are ambiguous
# (Excel's auto-detect could re-read "01-08-2022" as 1-Aug become a real
# date/mm-dd-yyyy), so those are pinned to the Text format first.
$safeRows  = @(3, 8, 9, 10, 11, 12, 17, 18, 19, 20, 21)
$ambigRows = @(4, 5, 6, 7, 13, 14, 15, 16)

foreach ($r in $safeRows) {
    $ws.Range("A$r").Replace("/", "-")
}

foreach ($r in $ambigRows) {
    $cell = $ws.Range("A$r")
    $cell.NumberFormat = "@"
    $cell.Replace("/", "-")
}

# --- Attendance flag updates (columns D, E, G, H) for the affected rows ---
$ws.Cells.Item(3, 4).Value  = 1   # D3
$ws.Cells.Item(3, 7).Value  = 1   # G3

$ws.Cells.Item(4, 4).Value  = 1   # D4
$ws.Cells.Item(4, 5).Value  = 1   # E4
$ws.Cells.Item(4, 8).Value  = 0   # H4

$ws.Cells.Item(5, 4).Value  = 1   # D5
$ws.Cells.Item(5, 5).Value  = 1   # E5
$ws.Cells.Item(5, 8).Value  = 0   # H5

$ws.Cells.Item(10, 4).Value = 1   # D10
$ws.Cells.Item(10, 5).Value = 1   # E10
$ws.Cells.Item(10, 8).Value = 0   # H10

$ws.Cells.Item(12, 4).Value = 1   # D12
$ws.Cells.Item(12, 7).Value = 1   # G12

$ws.Cells.Item(14, 4).Value = 1   # D14
$ws.Cells.Item(14, 5).Value = 1   # E14
$ws.Cells.Item(14, 8).Value = 0   # H14

$wb.Save()
